# Refresh the cryptocurrency price/volume table (columns B-E, rows 2-51)
# with the latest scraped values from coinranking.com.
# Column D (Price) values are forced to Text format ("@") before assignment
# so that numeric-looking strings (e.g. "8.20", "572.82") keep their exact
# textual representation instead of being auto-converted into numbers by
# Excel (which would silently drop meaningful trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '63.012.22'
$ws.Cells.Item(2, 5).Value = '  +5.49%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.419.50'
$ws.Cells.Item(3, 5).Value = '  +2.15%  '
$ws.Cells.Item(4, 5).Value = '  +0.55%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '572.82'
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '146.02'
$ws.Cells.Item(6, 5).Value = '  +6.16%  '
$ws.Cells.Item(7, 5).Value = '  -0.39%  '
$ws.Cells.Item(8, 5).Value = '  +2.49%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.452.85'
$ws.Cells.Item(9, 5).Value = '  +3.76%  '
$ws.Cells.Item(10, 5).Value = '  +5.91%  '
$ws.Cells.Item(11, 5).Value = '  +1.04%  '
$ws.Cells.Item(12, 5).Value = '  +3.01%  '
$ws.Cells.Item(13, 5).Value = '  +4.79%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '27.39'
$ws.Cells.Item(14, 5).Value = '  +7.21%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000178'
$ws.Cells.Item(15, 5).Value = '  +8.01%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '2.858.17'
$ws.Cells.Item(16, 5).Value = '  +2.28%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '62.854.66'
$ws.Cells.Item(17, 5).Value = '  +5.42%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.447.22'
$ws.Cells.Item(18, 5).Value = '  +3.95%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.91'
$ws.Cells.Item(19, 5).Value = '  -1.05%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '10.99'
$ws.Cells.Item(20, 5).Value = '  +5.19%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '328.55'
$ws.Cells.Item(21, 5).Value = '  +2.25%  '
$ws.Cells.Item(22, 5).Value = '  +2.32%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '2.04'
$ws.Cells.Item(23, 5).Value = '  +12.90%  '
$ws.Cells.Item(24, 5).Value = '  -0.34%  '
$ws.Cells.Item(25, 5).Value = '  +2.50%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '640.31'
$ws.Cells.Item(26, 5).Value = '  +14.54%  '
$ws.Cells.Item(27, 2).Value = 'Aptos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.53'
$ws.Cells.Item(27, 5).Value = '  +4.57%  '
$ws.Cells.Item(28, 2).Value = 'PEPE'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.0₃0988'
$ws.Cells.Item(28, 5).Value = '  +7.31%  '
$ws.Cells.Item(29, 2).Value = 'WrappedeETH'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.532.99'
$ws.Cells.Item(29, 5).Value = '  +2.00%  '
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '8.20'
$ws.Cells.Item(30, 5).Value = '  +2.64%  '
$ws.Cells.Item(31, 2).Value = 'Fetch.AI'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.41'
$ws.Cells.Item(31, 5).Value = '  +8.77%  '
$ws.Cells.Item(32, 2).Value = 'Kaspa'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.138'
$ws.Cells.Item(32, 5).Value = '  +6.07%  '
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.84'
$ws.Cells.Item(33, 5).Value = '  +3.71%  '
$ws.Cells.Item(34, 2).Value = 'ImmutableX'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.49'
$ws.Cells.Item(34, 5).Value = '  +5.11%  '
$ws.Cells.Item(35, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.995'
$ws.Cells.Item(35, 5).Value = '  -0.39%  '
$ws.Cells.Item(36, 2).Value = 'NEARProtocol'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '4.76'
$ws.Cells.Item(36, 5).Value = '  +5.20%  '
$ws.Cells.Item(37, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.374'
$ws.Cells.Item(37, 5).Value = '  +2.26%  '
$ws.Cells.Item(38, 2).Value = 'Monero'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '152.93'
$ws.Cells.Item(38, 5).Value = '  -0.02%  '
$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.41'
$ws.Cells.Item(39, 5).Value = '  +8.79%  '
$ws.Cells.Item(40, 2).Value = 'EthereumClassic'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '18.70'
$ws.Cells.Item(40, 5).Value = '  +3.03%  '
$ws.Cells.Item(41, 2).Value = 'dogwifhat'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.74'
$ws.Cells.Item(41, 5).Value = '  +14.16%  '
$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.77'
$ws.Cells.Item(42, 5).Value = '  +8.33%  '
$ws.Cells.Item(43, 2).Value = 'USDe'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  -0.04%  '
$ws.Cells.Item(44, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0₆0285'
$ws.Cells.Item(44, 5).Value = '  -4.42%  '
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '144.93'
$ws.Cells.Item(45, 5).Value = '  +4.63%  '
$ws.Cells.Item(46, 2).Value = 'Filecoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.60'
$ws.Cells.Item(46, 5).Value = '  +2.49%  '
$ws.Cells.Item(47, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '20.45'
$ws.Cells.Item(47, 5).Value = '  +7.59%  '
$ws.Cells.Item(48, 2).Value = 'Mantle'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.603'
$ws.Cells.Item(48, 5).Value = '  +3.28%  '
$ws.Cells.Item(49, 2).Value = 'Hedera'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0517'
$ws.Cells.Item(49, 5).Value = '  +3.47%  '
$ws.Cells.Item(50, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '12.68'
$ws.Cells.Item(50, 5).Value = '  +8.66%  '
$ws.Cells.Item(51, 2).Value = 'Stellar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0918'
$ws.Cells.Item(51, 5).Value = '  +2.51%  '
